# Applies the OOXML changes described by the commit:
#  - shift the "INSIGHTS" content placeholder on slide 2 down slightly
#  - clear the paragraph "space before" override on the affected text
#    boxes on slide 2 and slide 3

$p = $ppt.ActivePresentation

# --- Slide 2: "Content Placeholder 2" (shape id 16, the INSIGHTS textbox) ---
$s2 = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(3)

# Move the box down by 70340 EMU (≈5.538583 pt) -> off y 2672912 -> 2743252
$sh2.Top = $sh2.Top + (70340.0 / 12700.0)

# Remove the explicit "space before" on every paragraph in this text box
$tr2 = $sh2.TextFrame.TextRange
$count2 = $tr2.Paragraphs().Count
for ($i = 1; $i -le $count2; $i++) {
    $tr2.Paragraphs($i, 1).ParagraphFormat.SpaceBefore = 0
}
# the collection does not report the trailing empty paragraph - clear it too
$tr2.Paragraphs($count2 + 1, 1).ParagraphFormat.SpaceBefore = 0

# --- Slide 3: "Content Placeholder 2" (shape id 17, the bullet list textbox) ---
$s3 = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(6)

$tr3 = $sh3.TextFrame.TextRange
$count3 = $tr3.Paragraphs().Count
for ($i = 1; $i -le $count3; $i++) {
    $tr3.Paragraphs($i, 1).ParagraphFormat.SpaceBefore = 0
}
# same trailing-empty-paragraph quirk as above
$tr3.Paragraphs($count3 + 1, 1).ParagraphFormat.SpaceBefore = 0
